# Daily attendance processing - 2026-01-12 11:35:14
# Normalise the "Recorded By" column (G): entries that list the system
# account before the human editor's email are reordered so the human
# editor's email comes first, e.g.
#   "System, dnasr281@gmail.com"  ->  "dnasr281@gmail.com, System"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

$oldValue = "System, dnasr281@gmail.com"
$newValue = "dnasr281@gmail.com, System"

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    if ($cell.Text -eq $oldValue) {
        $cell.Value = $newValue
    }
}
